$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (15 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 11418.781
$ws.Cells.Item(62, 9).Value = 12965.962
$ws.Cells.Item(62, 11).Value = 12965.962
$ws.Cells.Item(62, 13).Value = -12341.962
$ws.Cells.Item(65, 8).Value = 11418.781
$ws.Cells.Item(65, 9).Value = 12965.962
$ws.Cells.Item(65, 11).Value = 64829.81
$ws.Cells.Item(65, 13).Value = -61709.81
$ws.Cells.Item(137, 8).Value = 29925.395
$ws.Cells.Item(137, 9).Value = 48564.094
$ws.Cells.Item(137, 10).Value = 6901.1177
$ws.Cells.Item(137, 11).Value = 145692.282
$ws.Cells.Item(137, 12).Value = 20703.3531
$ws.Cells.Item(137, 13).Value = -143142.282
$ws.Cells.Item(137, 14).Value = -25803.3531

# --- Sheet: ARM (21 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 58268.055
$ws.Cells.Item(74, 9).Value = 94028.27
$ws.Cells.Item(74, 10).Value = 2073.4285
$ws.Cells.Item(74, 11).Value = 94028.27
$ws.Cells.Item(74, 12).Value = 2073.4285
$ws.Cells.Item(74, 13).Value = -93154.27
$ws.Cells.Item(74, 14).Value = -3821.4285
$ws.Cells.Item(77, 8).Value = 58268.055
$ws.Cells.Item(77, 9).Value = 94028.27
$ws.Cells.Item(77, 10).Value = 2073.4285
$ws.Cells.Item(77, 11).Value = 470141.35
$ws.Cells.Item(77, 12).Value = 10367.1425
$ws.Cells.Item(77, 13).Value = -465773.35
$ws.Cells.Item(77, 14).Value = -19103.1425
$ws.Cells.Item(132, 8).Value = 2672834.5
$ws.Cells.Item(132, 9).Value = 3294872
$ws.Cells.Item(132, 10).Value = 919819.4399999999
$ws.Cells.Item(132, 11).Value = 9884616
$ws.Cells.Item(132, 12).Value = 2759458.32
$ws.Cells.Item(132, 13).Value = -9882086
$ws.Cells.Item(132, 14).Value = -2764518.32

# --- Sheet: BSM (34 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 2570
$ws.Cells.Item(22, 9).Value = 2922.8572
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 2922.8572
$ws.Cells.Item(22, 12).Value = 100
$ws.Cells.Item(22, 14).Value = -446
$ws.Cells.Item(47, 8).Value = 99800
$ws.Cells.Item(47, 10).Value = 99800
$ws.Cells.Item(47, 12).Value = 99800
$ws.Cells.Item(47, 14).Value = -100840
$ws.Cells.Item(61, 8).Value = 41141.4
$ws.Cells.Item(61, 10).Value = 41141.4
$ws.Cells.Item(61, 14).Value = -41767.4
$ws.Cells.Item(99, 8).Value = 888.1786
$ws.Cells.Item(99, 9).Value = 868.2174
$ws.Cells.Item(99, 10).Value = 980
$ws.Cells.Item(99, 11).Value = 868.2174
$ws.Cells.Item(99, 12).Value = 980
$ws.Cells.Item(99, 13).Value = 629.7826
$ws.Cells.Item(99, 14).Value = -3976
$ws.Cells.Item(107, 8).Value = 670.3077
$ws.Cells.Item(107, 9).Value = 720.1
$ws.Cells.Item(107, 10).Value = 504.33334
$ws.Cells.Item(107, 11).Value = 720.1
$ws.Cells.Item(107, 12).Value = 504.33334
$ws.Cells.Item(107, 13).Value = 1199.9
$ws.Cells.Item(107, 14).Value = -4344.33334
$ws.Cells.Item(134, 8).Value = 17667.746
$ws.Cells.Item(134, 9).Value = 1042.2982
$ws.Cells.Item(134, 10).Value = 112432.8
$ws.Cells.Item(134, 11).Value = 3126.8946
$ws.Cells.Item(134, 12).Value = 337298.4
$ws.Cells.Item(134, 13).Value = -591.8945999999996
$ws.Cells.Item(134, 14).Value = -342368.4

# --- Sheet: CRP (39 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8318.173000000001
$ws.Cells.Item(31, 9).Value = 6429.8223
$ws.Cells.Item(31, 10).Value = 14854.77
$ws.Cells.Item(31, 11).Value = 6429.8223
$ws.Cells.Item(31, 12).Value = 14854.77
$ws.Cells.Item(31, 13).Value = -6134.8223
$ws.Cells.Item(31, 14).Value = -15444.77
$ws.Cells.Item(34, 8).Value = 8318.173000000001
$ws.Cells.Item(34, 9).Value = 6429.8223
$ws.Cells.Item(34, 10).Value = 14854.77
$ws.Cells.Item(34, 11).Value = 6429.8223
$ws.Cells.Item(34, 12).Value = 14854.77
$ws.Cells.Item(34, 13).Value = -6227.8223
$ws.Cells.Item(34, 14).Value = -15258.77
$ws.Cells.Item(58, 8).Value = 1071.3492
$ws.Cells.Item(58, 9).Value = 706.0732
$ws.Cells.Item(58, 10).Value = 1752.091
$ws.Cells.Item(58, 11).Value = 706.0732
$ws.Cells.Item(58, 12).Value = 1752.091
$ws.Cells.Item(58, 13).Value = -503.0732
$ws.Cells.Item(58, 14).Value = -2158.091
$ws.Cells.Item(132, 8).Value = 1387.7805
$ws.Cells.Item(132, 9).Value = 943.9091
$ws.Cells.Item(132, 11).Value = 2831.7273
$ws.Cells.Item(132, 13).Value = -301.7273
$ws.Cells.Item(134, 8).Value = 1187.3829
$ws.Cells.Item(134, 9).Value = 1053.2972
$ws.Cells.Item(134, 10).Value = 1683.5
$ws.Cells.Item(134, 11).Value = 3159.8916
$ws.Cells.Item(134, 12).Value = 5050.5
$ws.Cells.Item(134, 13).Value = -624.8915999999999
$ws.Cells.Item(134, 14).Value = -10120.5
$ws.Cells.Item(136, 8).Value = 1071.3492
$ws.Cells.Item(136, 9).Value = 706.0732
$ws.Cells.Item(136, 10).Value = 1752.091
$ws.Cells.Item(136, 11).Value = 2118.2196
$ws.Cells.Item(136, 12).Value = 5256.272999999999
$ws.Cells.Item(136, 13).Value = 431.7803999999996
$ws.Cells.Item(136, 14).Value = -10356.273

# --- Sheet: CUL (64 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(123, 8).Value = 2553.158
$ws.Cells.Item(123, 9).Value = 1457.5
$ws.Cells.Item(123, 10).Value = 2845.3333
$ws.Cells.Item(123, 11).Value = 4372.5
$ws.Cells.Item(123, 12).Value = 8535.999899999999
$ws.Cells.Item(123, 13).Value = -1922.5
$ws.Cells.Item(123, 14).Value = -13435.9999
$ws.Cells.Item(129, 8).Value = 28154.756
$ws.Cells.Item(129, 10).Value = 34585.535
$ws.Cells.Item(129, 12).Value = 103756.605
$ws.Cells.Item(129, 14).Value = -113756.605
$ws.Cells.Item(131, 8).Value = 32052164
$ws.Cells.Item(131, 9).Value = 506.66666
$ws.Cells.Item(131, 10).Value = 34723136
$ws.Cells.Item(131, 11).Value = 1519.99998
$ws.Cells.Item(131, 12).Value = 104169408
$ws.Cells.Item(131, 13).Value = 3520.00002
$ws.Cells.Item(131, 14).Value = -104179488
$ws.Cells.Item(133, 8).Value = 2845.4707
$ws.Cells.Item(133, 9).Value = 1444.0769
$ws.Cells.Item(133, 10).Value = 7400
$ws.Cells.Item(133, 11).Value = 4332.2307
$ws.Cells.Item(133, 12).Value = 22200
$ws.Cells.Item(133, 13).Value = 727.7692999999999
$ws.Cells.Item(133, 14).Value = -32320
$ws.Cells.Item(134, 8).Value = 4298.775
$ws.Cells.Item(134, 9).Value = 1527.5416
$ws.Cells.Item(134, 11).Value = 4582.6248
$ws.Cells.Item(134, 13).Value = 487.3752000000004
$ws.Cells.Item(137, 8).Value = 8689441
$ws.Cells.Item(137, 9).Value = 33334234
$ws.Cells.Item(137, 10).Value = 4125590.5
$ws.Cells.Item(137, 11).Value = 100002702
$ws.Cells.Item(137, 12).Value = 12376771.5
$ws.Cells.Item(137, 13).Value = -99997602
$ws.Cells.Item(137, 14).Value = -12386971.5
$ws.Cells.Item(138, 8).Value = 14495617
$ws.Cells.Item(138, 9).Value = 1865
$ws.Cells.Item(138, 10).Value = 22225618
$ws.Cells.Item(138, 11).Value = 5595
$ws.Cells.Item(138, 12).Value = 66676854
$ws.Cells.Item(138, 13).Value = -455
$ws.Cells.Item(138, 14).Value = -66687134
$ws.Cells.Item(139, 8).Value = 20373418
$ws.Cells.Item(139, 9).Value = 33335554
$ws.Cells.Item(139, 10).Value = 7411282.5
$ws.Cells.Item(139, 11).Value = 100006662
$ws.Cells.Item(139, 12).Value = 22233847.5
$ws.Cells.Item(139, 13).Value = -100001522
$ws.Cells.Item(139, 14).Value = -22244127.5
$ws.Cells.Item(140, 8).Value = 15630318
$ws.Cells.Item(140, 9).Value = 3061.25
$ws.Cells.Item(140, 10).Value = 20839404
$ws.Cells.Item(140, 11).Value = 9183.75
$ws.Cells.Item(140, 12).Value = 62518212
$ws.Cells.Item(140, 13).Value = -4003.75
$ws.Cells.Item(140, 14).Value = -62528572
$ws.Cells.Item(141, 8).Value = 4906650
$ws.Cells.Item(141, 9).Value = 1002.6
$ws.Cells.Item(141, 10).Value = 6950669.5
$ws.Cells.Item(141, 11).Value = 3007.8
$ws.Cells.Item(141, 12).Value = 20852008.5
$ws.Cells.Item(141, 13).Value = 2172.2
$ws.Cells.Item(141, 14).Value = -20862368.5

# --- Sheet: GSM (18 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 3707.1082
$ws.Cells.Item(70, 9).Value = 3671.4194
$ws.Cells.Item(70, 10).Value = 3891.5
$ws.Cells.Item(70, 11).Value = 3671.4194
$ws.Cells.Item(70, 12).Value = 3891.5
$ws.Cells.Item(70, 13).Value = -3401.4194
$ws.Cells.Item(70, 14).Value = -4431.5
$ws.Cells.Item(73, 8).Value = 3707.1082
$ws.Cells.Item(73, 9).Value = 3671.4194
$ws.Cells.Item(73, 10).Value = 3891.5
$ws.Cells.Item(73, 11).Value = 3671.4194
$ws.Cells.Item(73, 12).Value = 3891.5
$ws.Cells.Item(73, 13).Value = -2735.4194
$ws.Cells.Item(73, 14).Value = -5763.5
$ws.Cells.Item(102, 8).Value = 14879.031
$ws.Cells.Item(102, 9).Value = 7639.048
$ws.Cells.Item(102, 11).Value = 7639.048
$ws.Cells.Item(102, 13).Value = -6017.048

# --- Sheet: LTW (24 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(60, 8).Value = 21533.334
$ws.Cells.Item(60, 10).Value = 21533.334
$ws.Cells.Item(60, 14).Value = -22551.334
$ws.Cells.Item(100, 8).Value = 29391.73
$ws.Cells.Item(100, 9).Value = 113032.78
$ws.Cells.Item(100, 10).Value = 2507.1072
$ws.Cells.Item(100, 11).Value = 113032.78
$ws.Cells.Item(100, 12).Value = 2507.1072
$ws.Cells.Item(100, 13).Value = -112491.78
$ws.Cells.Item(100, 14).Value = -3589.1072
$ws.Cells.Item(132, 8).Value = 199541.22
$ws.Cells.Item(132, 9).Value = 47288
$ws.Cells.Item(132, 10).Value = 593608.4
$ws.Cells.Item(132, 11).Value = 141864
$ws.Cells.Item(132, 12).Value = 1780825.2
$ws.Cells.Item(132, 13).Value = -139334
$ws.Cells.Item(132, 14).Value = -1785885.2
$ws.Cells.Item(136, 8).Value = 140053.31
$ws.Cells.Item(136, 9).Value = 164865.36
$ws.Cells.Item(136, 10).Value = 2459.182
$ws.Cells.Item(136, 11).Value = 494596.08
$ws.Cells.Item(136, 12).Value = 7377.545999999999
$ws.Cells.Item(136, 13).Value = -492046.08
$ws.Cells.Item(136, 14).Value = -12477.546

# --- Sheet: WVR (29 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(59, 8).Value = 12500
$ws.Cells.Item(59, 10).Value = 12500
$ws.Cells.Item(59, 12).Value = 12500
$ws.Cells.Item(59, 14).Value = -13976
$ws.Cells.Item(61, 8).Value = 7790
$ws.Cells.Item(61, 10).Value = 7790
$ws.Cells.Item(61, 12).Value = 7790
$ws.Cells.Item(61, 14).Value = -8374
$ws.Cells.Item(126, 8).Value = 962
$ws.Cells.Item(126, 9).Value = 703.3333
$ws.Cells.Item(126, 10).Value = 1350
$ws.Cells.Item(126, 11).Value = 2109.9999
$ws.Cells.Item(126, 12).Value = 4050
$ws.Cells.Item(126, 13).Value = 360.0001000000002
$ws.Cells.Item(126, 14).Value = -8990
$ws.Cells.Item(132, 8).Value = 3312.1365
$ws.Cells.Item(132, 9).Value = 715.0357
$ws.Cells.Item(132, 10).Value = 7857.0625
$ws.Cells.Item(132, 11).Value = 2145.1071
$ws.Cells.Item(132, 12).Value = 23571.1875
$ws.Cells.Item(132, 13).Value = 384.8928999999998
$ws.Cells.Item(132, 14).Value = -28631.1875
$ws.Cells.Item(136, 8).Value = 1509112.8
$ws.Cells.Item(136, 9).Value = 1880966.2
$ws.Cells.Item(136, 10).Value = 625960.9
$ws.Cells.Item(136, 11).Value = 5642898.6
$ws.Cells.Item(136, 12).Value = 1877882.7
$ws.Cells.Item(136, 13).Value = -5640348.6
$ws.Cells.Item(136, 14).Value = -1882982.7
